$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The B4 cell (FilesTab query) had its query text corrected: the
# `File Type` and `Breed` output columns were removed from the RETURN
# clause (ICDC Breed script correction).
$newFileQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Bouvier des Flandres'] `nOPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN  coalesce(f.file_name, '') AS ``File Name``,`n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newFileQuery

# The row height of row 4 shrinks because the corrected query text wraps
# onto fewer lines.
$ws.Rows.Item(4).RowHeight = 217.5

# Update the view: scroll so row 4 is at the top and select B4 (matches
# the saved sheet view state after the edit).
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
